$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting existing rows 7-16 down to 8-17
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the invoice number field metadata
$ws.Range("A7").Value = "Číslo faktury"
$ws.Range("C7").Value = "Invoice Number"
$ws.Range("C7").Style = $ws.Range("C6").Style
$ws.Range("B7").Value = "Invoice_number"

# Update selection to match the post-edit state (D17)
$ws.Range("D17").Select()
